$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '62.396.12'
$ws.Range('E2').Value = '  +0.50%  '
$ws.Range('D3').Value = '3.016.14'
$ws.Range('E3').Value = '  +0.64%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '0.999'
$ws.Range('E4').Value = '  -0.25%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '598.41'
$ws.Range('E5').Value = '  +2.95%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '147.49'
$ws.Range('E6').Value = '  +1.73%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.999'
$ws.Range('E7').Value = '  -0.13%  '
$ws.Range('B8').Value = 'LidoStakedEther'
$ws.Range('C8').Value = 'https://coinranking.com/coin/VINVMYf0u+lidostakedether-steth'
$ws.Range('D8').Value = '3.014.52'
$ws.Range('E8').Value = '  +0.64%  '
$ws.Range('B9').Value = 'XRP'
$ws.Range('C9').Value = 'https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp'
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.518'
$ws.Range('E9').Value = '  -1.60%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '6.31'
$ws.Range('E10').Value = '  +9.46%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.150'
$ws.Range('E11').Value = '  +1.72%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.457'
$ws.Range('E12').Value = '  -0.67%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.0000234'
$ws.Range('E13').Value = '  +2.98%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '34.57'
$ws.Range('E14').Value = '  +0.77%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.126'
$ws.Range('E15').Value = '  +2.64%  '
$ws.Range('D16').Value = '3.507.67'
$ws.Range('E16').Value = '  +0.46%  '
$ws.Range('B17').Value = 'Polkadot'
$ws.Range('C17').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '7.00'
$ws.Range('E17').Value = '  -1.35%  '
$ws.Range('B18').Value = 'WrappedBTC'
$ws.Range('C18').Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range('D18').Value = '62.242.59'
$ws.Range('E18').Value = '  +0.18%  '
$ws.Range('D19').Value = '3.009.23'
$ws.Range('E19').Value = '  +0.20%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '451.19'
$ws.Range('E20').Value = '  -1.46%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '14.08'
$ws.Range('E21').Value = '  +1.17%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '0.688'
$ws.Range('E22').Value = '  +0.43%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '7.40'
$ws.Range('E23').Value = '  -0.24%  '
$ws.Range('B24').Value = 'RenderToken'
$ws.Range('C24').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '11.32'
$ws.Range('E24').Value = '  +13.24%  '
$ws.Range('B25').Value = 'Litecoin'
$ws.Range('C25').Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '81.95'
$ws.Range('E25').Value = '  +0.53%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '2.29'
$ws.Range('E26').Value = '  +4.12%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '12.28'
$ws.Range('E27').Value = '  +0.18%  '
$ws.Range('E28').Value = '  +0.16%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '2.72'
$ws.Range('E29').Value = '  +4.43%  '
$ws.Range('B30').Value = 'NEARProtocol'
$ws.Range('C30').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '7.32'
$ws.Range('E30').Value = '  +5.42%  '
$ws.Range('B31').Value = 'FirstDigitalUSD'
$ws.Range('C31').Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '0.998'
$ws.Range('E31').Value = '  -0.32%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '2.10'
$ws.Range('E32').Value = '  +1.97%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '27.44'
$ws.Range('E33').Value = '  -2.99%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.111'
$ws.Range('E34').Value = '  +3.27%  '
$ws.Range('D35').Value = '0.0₃0860'
$ws.Range('E35').Value = '  +8.97%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '1.03'
$ws.Range('E36').Value = '  +0.69%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '5.85'
$ws.Range('E37').Value = '  +2.15%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '50.48'
$ws.Range('E38').Value = '  +0.53%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '2.07'
$ws.Range('E39').Value = '  -1.31%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '9.02'
$ws.Range('E40').Value = '  -1.67%  '
$ws.Range('B41').Value = 'dogwifhat'
$ws.Range('C41').Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '2.97'
$ws.Range('E41').Value = '  +4.14%  '
$ws.Range('B42').Value = 'Kaspa'
$ws.Range('C42').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.124'
$ws.Range('E42').Value = '  +8.04%  '
$ws.Range('B43').Value = 'Arweave'
$ws.Range('C43').Value = 'https://coinranking.com/coin/7XWg41D1+arweave-ar'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '41.66'
$ws.Range('E43').Value = '  +14.90%  '
$ws.Range('B44').Value = 'Bittensor'
$ws.Range('C44').Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '403.72'
$ws.Range('E44').Value = '  +3.77%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.278'
$ws.Range('E45').Value = '  +3.39%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.0353'
$ws.Range('E46').Value = '  -0.78%  '
$ws.Range('D47').Value = '2.723.22'
$ws.Range('E47').Value = '  +0.20%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '132.05'
$ws.Range('E48').Value = '  +3.12%  '
$ws.Range('E49').Value = '  +0.11%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '2.20'
$ws.Range('E50').Value = '  +0.17%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.108'
$ws.Range('E51').Value = '  -1.04%  '
